$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Build the new "Names" column (B) and insert a new row (old row count 7 -> 8) ---
# Target layout (row -> A value, B value, hidden?):
#  1 Numbers  Names    (header, not hidden)
#  2 5        Jacques  (not hidden)
#  3 5        Alex     (hidden)
#  4 4        Patrick  (hidden)
#  5 3        Jack     (not hidden)
#  6 3        Neil     (hidden)
#  7 2        John     (hidden)
#  8 1        Jason    (not hidden)

$ws.Range("B1").Value = "Names"
$ws.Range("B2").Value = "Jacques"
$ws.Range("B3").Value = "Alex"
$ws.Range("B4").Value = "Patrick"
$ws.Range("B5").Value = "Jack"
$ws.Range("B6").Value = "Neil"
$ws.Range("B7").Value = "John"
$ws.Range("B8").Value = "Jason"

$ws.Range("A8").Value = $ws.Range("A7").Value()
$ws.Range("A7").Value = $ws.Range("A6").Value()
$ws.Range("A6").Value = $ws.Range("A5").Value()
$ws.Range("A5").Value = $ws.Range("A4").Value()
$ws.Range("A4").Value = $ws.Range("A3").Value()
$ws.Range("A3").Value = $ws.Range("A2").Value()

# --- Reset and re-apply the autofilter over the new A1:B8 range ---
$r = $ws.Range("A1:B8")
$r.AutoFilter()
$r.AutoFilter(1, @("3", "1", "5"))
$r.AutoFilter(2, "J*")

# --- Fix up row hidden states to match the target exactly ---
$ws.Rows.Item(2).Hidden = $false
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(4).Hidden = $true
$ws.Rows.Item(5).Hidden = $false
$ws.Rows.Item(6).Hidden = $true
$ws.Rows.Item(7).Hidden = $true
$ws.Rows.Item(8).Hidden = $false

# --- Update the workbook-level _xlnm._FilterDatabase defined name for this sheet ---
for ($i = 1; $i -le $wb.Names.Count; $i++) {
  $n = $wb.Names.Item($i)
  if ($n.Name() -eq "Single Column Numbers!_FilterDatabase") {
    $n.RefersTo = "='Single Column Numbers'!`$A`$1:`$B`$8"
  }
}
